$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.176.98"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.601.67"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'211.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "'0.0614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'18.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "1.823.42"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.600.80"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "'0.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "26.160.98"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'61.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'203.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.82%  "
$ws.Range("D25").Value = "'144.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -7.34%  "
$ws.Range("D28").Value = "'15.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'3.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "1.138.95"
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  +6.11%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'0.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").Value = "'0.783"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "1.738.13"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "'92.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'0.0506"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'0.407"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0949"
$ws.Range("E51").Value = "  -15.60%  "
